# Updates crypto price/volume figures per the Tue Jan 17 15:55:40 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "303.35"
    "E2" = "2.08%"
    "D3" = "32.06"
    "E3" = "2.33%"
    "D4" = "5.110"
    "E4" = "0.47%"
    "D5" = "0.07824"
    "E5" = "-2.46%"
    "D6" = "2.264"
    "E6" = "-14.23%"
    "D7" = "7.834"
    "E7" = "0.45%"
    "E8" = "0.55%"
    "D9" = "0.9247"
    "E9" = "-0.51%"
    "D10" = "0.1762"
    "E10" = "1.13%"
    "D11" = "0.07748"
    "E11" = "7.99%"
    "D12" = "0.08858"
    "E12" = "-1.11%"
    "D13" = "0.03100"
    "E13" = "0.56%"
    "D14" = "0.09995"
    "E14" = "-0.16%"
    "D15" = "0.001513"
    "E15" = "0.65%"
    "D16" = "0.006008"
    "E16" = "-0.17%"
    "D17" = "3.451"
    "E17" = "-2.65%"
    "D18" = "2.244"
    "E18" = "-0.13%"
    "D19" = "0.3273"
    "E19" = "1.30%"
    "D20" = "0.1339"
    "E20" = "-0.26%"
    "E21" = "8.66%"
    "D22" = "0.1821"
    "E22" = "12.43%"
    "D23" = "0.04595"
    "E23" = "0.30%"
    "D24" = "0.001253"
    "E24" = "1.00%"
    "D25" = "0.004486"
    "E25" = "1.48%"
    "D26" = "0.0001252"
    "E26" = "4.45%"
    "E27" = "-1.00%"
    "D39" = "0.01788"
    "E39" = "0.80%"
    "D40" = "0.04809"
    "E40" = "7.37%"
    "D41" = "0.007189"
    "E41" = "4.83%"
    "D42" = "0.1367"
    "E42" = "1.76%"
    "D43" = "0.002124"
    "E43" = "-0.43%"
    "D44" = "0.009984"
    "E44" = "4.07%"
    "D45" = "0.00006281"
    "E45" = "-3.75%"
    "D46" = "0.00000000751"
    "E46" = "0.46%"
    "D47" = "0.003601"
    "E47" = "-58.71%"
    "D48" = "1.162"
    "E48" = "41.58%"
    "D49" = "0.00002103"
    "E49" = "0.46%"
    "D50" = "0.0002003"
    "E50" = "0.46%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
